$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.394.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.567.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.93%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("E6").Value = "  -2.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3671"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.54"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3370"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.46%  "
$ws.Range("E10").Value = "  -4.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07560"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.74%  "
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.045"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.853"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001142"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.573.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -8.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06697"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.246"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.5272"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -9.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "22.409.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.397"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.918"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.80%  "
$ws.Range("E29").Value = "  -4.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.950"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "124.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.752.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.236"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -10.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.981"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9832"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.31"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -12.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08409"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02529"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2293"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.96%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06505"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.34%  "
$ws.Range("B41").Value = "InternetComputer(DFINITY)"
$ws.Range("C41").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.505"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.241"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6370"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6018"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.771"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.112"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "121.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07266"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.04%  "
